# The deck currently carries two theme parts: the "Integral" theme that is
# actually applied to the slide master/design (ppt/theme/theme2.xml) and an
# unused stock "Office Theme" colour set (ppt/theme/theme1.xml). The edit
# being replayed swaps which theme file holds which palette, so that the
# design that is applied to the presentation ends up using the plain
# "Office Theme" palette instead of "Integral".
#
# The only thing that actually changes inside the package is the 12-colour
# theme colour scheme of the applied design (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) -- everything else in both theme parts (font scheme,
# format scheme, names) is byte-identical between the two files. So we
# recolor the live design's theme color scheme to the "Office Theme" values.
#
# NOTE: PowerPoint's ColorFormat.RGB is a COM "OLE_COLOR" / 0xBBGGRR value
# (blue/green/red byte order), not 0xRRGGBB.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$colors = $design.SlideMaster.Theme.ThemeColorScheme

# Office Theme palette (target), expressed as 0xBBGGRR for the RGB setter.
$colors.Item(1).RGB  = 0x000000   # dk1      -> #000000
$colors.Item(2).RGB  = 0xFFFFFF   # lt1      -> #FFFFFF
$colors.Item(3).RGB  = 0x6A5444   # dk2      -> #44546A
$colors.Item(4).RGB  = 0xE6E6E7   # lt2      -> #E7E6E6
$colors.Item(5).RGB  = 0xD59B5B   # accent1  -> #5B9BD5
$colors.Item(6).RGB  = 0x317DED   # accent2  -> #ED7D31
$colors.Item(7).RGB  = 0xA5A5A5   # accent3  -> #A5A5A5
$colors.Item(8).RGB  = 0x00C0FF   # accent4  -> #FFC000
$colors.Item(9).RGB  = 0xC47244   # accent5  -> #4472C4
$colors.Item(10).RGB = 0x47AD70   # accent6  -> #70AD47
$colors.Item(11).RGB = 0xC16305   # hlink    -> #0563C1
$colors.Item(12).RGB = 0x724F95   # folHlink -> #954F72
